# Apply the "Trade #15 closed" update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.27   # Total P&L %
$wsSummary.Range("B6").Value = 15      # Total Trades
$wsSummary.Range("B9").Value = 26.67   # Win Rate %

# --- Strategy Status sheet ------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 15       # Trades for MarketMaking
$wsStatus.Range("G4").Value = 26.67    # Win Rate % for MarketMaking

# --- Helper to append the new trade row (#15) -----------------------------
function Add-Trade15Row($ws) {
    $ws.Cells.Item(16, 1).Value = 15

    # The Date column holds a plain text value (e.g. "2026-02-17"), not a
    # real Excel date. Force text format first so Excel's autodetection
    # doesn't silently convert the string into a date serial number, then
    # restore the default style so no stray formatting is left behind.
    $dateCell = $ws.Cells.Item(16, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item(16, 3).Value = "15:17:43"
    $ws.Cells.Item(16, 4).Value = "MarketMaking"
    $ws.Cells.Item(16, 5).Value = "UP"
    $ws.Cells.Item(16, 6).Value = 0.88
    $ws.Cells.Item(16, 7).Value = 0.88
    $ws.Cells.Item(16, 8).Value = "CLOSED"
    $ws.Cells.Item(16, 9).Value = 0
    $ws.Cells.Item(16, 10).Value = 0
    $ws.Cells.Item(16, 11).Value = 99.8
    $ws.Cells.Item(16, 12).Value = 0
    $ws.Cells.Item(16, 13).Value = 0
    $ws.Cells.Item(16, 14).Value = 0.6
    $ws.Cells.Item(16, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(16, 16).Value = "early_exit"
    $ws.Cells.Item(16, 17).Value = 0.14
}

# --- All Trades sheet -------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade15Row $wsAllTrades

# --- MarketMaking sheet -------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade15Row $wsMarketMaking
